$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(99, 1).Value = "2025-04-29 16:57:42"
$ws.Cells.Item(99, 2).Value = 259
